# Update the "Förändrad" (Changed) date column (C) for rows 2-27
# from 45335 (2024-02-13) to 45336 (2024-02-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45335) {
        $cell.Value2 = 45336
    }
}
